$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("findNewCarTest"): reorder brand rows and add a new "tata" row ---
$ws1.Range("C3").Value = "toyota"
$ws1.Range("D3").Value = "Toyota Cars"
$ws1.Range("C4").Value = "kia"
$ws1.Range("D4").Value = "Kia Cars"
$ws1.Range("C5").Value = "honda"
$ws1.Range("D5").Value = "Honda Cars"
$ws1.Range("A6").Value = "chrome"
$ws1.Range("B6").Value = "Y"
$ws1.Range("C6").Value = "tata"
$ws1.Range("D6").Value = "Tata Cars"

# --- Sheet2 ("carNameAndPrice"): flip runmode to Y and add a new "tata" row ---
$ws2.Range("B2").Value = "Y"
$ws2.Range("B3").Value = "Y"
$ws2.Range("B4").Value = "Y"
$ws2.Range("B5").Value = "Y"
$ws2.Range("A6").Value = "chrome"
$ws2.Range("B6").Value = "Y"
$ws2.Range("C6").Value = "tata"

# --- View state: switch the active/selected tab from sheet1 to sheet2 ---
$ws1.Range("S1:X9").Select()
$ws2.Activate()
$ws2.Range("C7").Select()
